# "Avance javascript y se agregó link"
# Adds a new row (21) with a JavaScript "switch" reference link + description,
# right below the existing list of links in column B/C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://developer.mozilla.org/es/docs/Web/JavaScript/Reference/Statements/switch"
$newDesc = "Ejemplos de switch en javascript"

# Write the new values first (B21 = url, C21 = description).
$ws.Range("B21").Value = $newUrl
$ws.Range("C21").Value = $newDesc

# Turn B21 into a real hyperlink, same as every other link in column B.
$ws.Hyperlinks.Add($ws.Range("B21"), $newUrl)

# Match the formatting used by the rest of the table: hyperlink style for
# column B, plain style for column C (re-apply after Hyperlinks.Add, which
# otherwise stamps its own style onto B21).
$ws.Range("B21").Style = $ws.Range("B20").Style
$ws.Range("C21").Style = $ws.Range("C20").Style

# Match the selection left behind in the saved file.
$ws.Range("C20").Select()
